# Update ticket/interest counts and mark a sold-out show as "不可售"
# across the relevant worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1237
$ws1.Range("F5").Value  = 1063
$ws1.Range("F6").Value  = 14018
$ws1.Range("F7").Value  = 15327
$ws1.Range("F11").Value = 183
$ws1.Range("F18").Value = 72
$ws1.Range("F19").Value = 25
$ws1.Range("F20").Value = 1183
$ws1.Range("F23").Value = 5957
$ws1.Range("F26").Value = 5510
$ws1.Range("F28").Value = 136
$ws1.Range("F29").Value = 92
$ws1.Range("F30").Value = 434

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value  = "不可售"
$ws4.Range("F5").Value  = 1237
$ws4.Range("F6").Value  = 1063
$ws4.Range("F7").Value  = 14018
$ws4.Range("F8").Value  = 15327
$ws4.Range("F12").Value = 183
$ws4.Range("F19").Value = 72
$ws4.Range("F20").Value = 25
$ws4.Range("F21").Value = 1183
$ws4.Range("F25").Value = 5957
$ws4.Range("F28").Value = 5510
$ws4.Range("F30").Value = 136
$ws4.Range("F31").Value = 92
$ws4.Range("F32").Value = 434
